# This script applies a cyclic re-shuffle of the weekly price rows
# (rows 2-9) in the "Fruta, Comercializadora del Agro de Limarí - Ciruela"
# sheet. The identifying columns (Mercado, Región, Codreg, Tipo, Producto,
# Categoría, Variedad, Origen) are identical across all rows, only the
# per-record fields (Fecha, Calidad, Volumen, Precio mínimo/máximo/
# promedio ponderado, Unidad de comercialización, Precio $/Kg, Kg/unidad)
# move between rows, per the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values of the columns that move, for rows 2..9.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")
$before = @{}
for ($r = 2; $r -le 9; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $row
}

# Destination row -> source row (cyclic permutation observed in the diff).
$mapping = @{
    2 = 5
    3 = 6
    4 = 7
    5 = 2
    6 = 3
    7 = 8
    8 = 9
    9 = 4
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $before[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
